$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Item", "Price"),
    @("Chicken", 6.99),
    @("Shrimp", 7.49),
    @("Beef & shrimp", 8.29),
    @("Beef", 7.19),
    @("Chicken & Shrimp", 8.19),
    @("Spring roll (1)", 0.99),
    @("Chicken egg roll (1)", 1.29),
    @("Chicken dumplings (6)", 2.49),
    @("Shrimp Tempura (3)", 2.99)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(2).ColumnWidth = 15

$ws.Range("B11").Select() | Out-Null

$win = $wb.Windows.Item(1)
$win.Left = 7215
$win.Top = 3600
